$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.074.46"

$ws.Range("D3").Value = "1.971.51"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.75"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4953"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4207"
$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("E9").Value = "  +3.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09279"
$ws.Range("E10").Value = "  +5.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.098"
$ws.Range("E11").Value = "  -1.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.84"
$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("D13").Value = "1.989.98"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("E14").Value = "  -1.46%  "

$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.012"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("E18").Value = "  -4.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06719"
$ws.Range("E19").Value = "  +1.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.19"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.961"
$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").Value = "29.108.76"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.264"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").Value = "2.207.18"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("E27").Value = "  +1.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.72"
$ws.Range("E28").Value = "  -1.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.260"
$ws.Range("E29").Value = "  -4.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.264"
$ws.Range("E30").Value = "  -3.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.14"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("E32").Value = "  -0.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09853"
$ws.Range("E33").Value = "  -0.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.508"
$ws.Range("E34").Value = "  -3.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.812"
$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.732"
$ws.Range("E36").Value = "  -1.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02422"
$ws.Range("E37").Value = "  -0.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.329"
$ws.Range("E38").Value = "  +3.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06427"
$ws.Range("E39").Value = "  +1.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.032"
$ws.Range("E40").Value = "  -5.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6477"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.50"
$ws.Range("E42").Value = "  -1.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2005"
$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.010"
$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.365"
$ws.Range("E45").Value = "  +7.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6208"
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("E48").Value = "  -0.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.481"
$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000327"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06975"
$ws.Range("E51").Value = "  -0.12%  "
